# Commit: "Fix: adding ENUM type on Konten"
# Adds a "type" column (header in E1, values "Film" / "Serial_TV" for each
# content row) to the Konten sheet, replacing the old "sutradara" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Konten")

# Header: sutradara -> type
$ws.Range("E1").Value = "type"

# Film rows: 2-134 and 288-310
$ws.Range("E2:E134").Value = "Film"
$ws.Range("E288:E310").Value = "Film"

# Serial_TV rows: 135-287
$ws.Range("E135:E287").Value = "Serial_TV"

Write-Host "Konten E column updated"
